$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("4.2"): combine the two task descriptions into one cell,
# reordering "Xay dung module quan ly nguoi dung" before the UI-design text.
$ws.Range("D4").Value = 'Xây dựng module quản lý người dùng.    Thiết kế giao diện cửa sổ chính (màn hình đăng nhập, các menu trỏ đến các module con,…). '

# Row 9 ("4.7"): assign to both Huyen and Do instead of just Lan.
$ws.Range("F9").Value = 'Huyền + Độ '

# Row 5 grew taller once the content moved off of it (its text wraps
# differently now the neighboring rows reflowed).
$ws.Rows.Item(5).RowHeight = 66.75

# Update the view: scroll position resets and selection moves to F10.
$ws.Range("F10").Select()
